# BreakoutBoard_RevB_BOM.xlsx edit:
#  - drop the two unused blank sheets
#  - rename the BOM sheet for the Breakout Board Rev B ("Smoky")
#  - retitle the sheet's heading cell to match
#  - swap the J1/J2 connector part data (rows 4 & 5) - the TOP/BOT
#    connector details had been entered against the wrong reference
#    designators

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the empty placeholder sheets.
[void]$wb.Worksheets("Sheet2").Delete()
[void]$wb.Worksheets("Sheet3").Delete()

# Rename the remaining (BOM) sheet.
$ws = $wb.Worksheets("Power Supply Board - Rev A")
$ws.Name = "Breakout Board - Rev B"

# Update the big merged title cell (A1:L1).
$ws.Range("A1").Value = "Bill of Materials for 'Marmote - Breakout Board Rev B (Smoky)'"

# Rows 4 and 5 (J1 / J2) had their Value/Package, Manufacturer Part #,
# Supplier Part #, Details and Unit Price swapped between the two
# connectors - fix by exchanging those columns between the two rows.
$swapCols = @("C", "D", "F", "H", "I", "K")
foreach ($col in $swapCols) {
    $cellTop = $ws.Range($col + "4")
    $cellBottom = $ws.Range($col + "5")
    $valueTop = $cellTop.Value2
    $valueBottom = $cellBottom.Value2
    $cellTop.Value = $valueBottom
    $cellBottom.Value = $valueTop
}

# Leave the title range selected, matching the post-edit view state.
[void]$ws.Range("A1:L1").Select()
